$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.261.61"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "'2.295.57"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'507.38"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'129.68"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "'2.318.79"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").Value = "'0.0981"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "'5.09"
$ws.Range("E12").Value = "  +8.13%  "
$ws.Range("D13").Value = "'0.341"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "'23.71"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("D15").Value = "'2.700.67"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "'55.287.65"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "'2.305.63"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'10.47"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").Value = "'4.18"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'312.80"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").Value = "'6.65"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("D24").Value = "'60.27"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "'0.994"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +3.47%  "
$ws.Range("D27").Value = "'7.52"
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("D28").Value = "'172.53"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "'0.0₃0713"
$ws.Range("E29").Value = "  +4.21%  "
$ws.Range("B30").Value = "'Aptos"
$ws.Range("C30").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'6.16"
$ws.Range("E30").Value = "  +4.22%  "
$ws.Range("B31").Value = "'Fetch.AI"
$ws.Range("C31").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +7.38%  "
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D34").Value = "'18.04"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "'1.24"
$ws.Range("E36").Value = "  +4.10%  "
$ws.Range("D37").Value = "'0.920"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").Value = "'3.90"
$ws.Range("E38").Value = "  +5.47%  "
$ws.Range("D39").Value = "'36.94"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("D41").Value = "'0.377"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "'137.05"
$ws.Range("E42").Value = "  +9.50%  "
$ws.Range("D43").Value = "'5.15"
$ws.Range("E43").Value = "  +6.69%  "
$ws.Range("D44").Value = "'3.47"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "'262.39"
$ws.Range("E45").Value = "  +10.08%  "
$ws.Range("D46").Value = "'0.0509"
$ws.Range("E46").Value = "  +3.40%  "
$ws.Range("D47").Value = "'0.0917"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "'0.554"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "'0.376"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Value = "'0.0212"
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("D51").Value = "'16.62"
$ws.Range("E51").Value = "  +2.61%  "
